$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$samplesQuery = @'
MATCH (ss:study_subject)
WITH COLLECT(ss.study_subject_id) AS all_subjects
MATCH (samp:sample)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)-[:study_of_program]->(p)
MATCH (samp)<-[:file_of_sample]-(f)-[:file_of_laboratory_procedure]->(lp)
MATCH (ss)<-[:diagnosis_of_study_subject]-(d)
MATCH (d)<-[:tp_of_diagnosis]-(tp)
WHERE tp.chemotherapy_regimen IN ["FEC (3 week cycles)"]
    distinct lp,
    toInteger(split(ss.study_subject_id,'-')[2]) AS subject_id_num,
    collect(distinct f.file_id) AS files,
    samp, ss, s, p, all_subjects
RETURN
 samp.sample_id AS `Sample ID`,
            ss.study_subject_id AS `Case ID`,
            p.program_acronym AS `Program Code`,
            s.study_acronym AS `Arm`,
            ss.disease_subtype AS `Diagnosis`,
            samp.tissue_type AS `Tissue Type`,
            samp.composition AS `Tissue Composition`,
            samp.sample_anatomic_site AS `Sample Anatomic Site`,
            samp.method_of_sample_procurement AS `Sample Procurement Method`
'@

$filesQuery = @'
MATCH (f:file)-->(parent)
MATCH (f)-[:file_of_sample]->(samp)
MATCH (samp)-[:sample_of_study_subject]->(ss)
MATCH (ss)-[:study_subject_of_study]->(s)
MATCH (s)-[:study_of_program]->(p)
MATCH (d)-[:diagnosis_of_study_subject]->(ss)
MATCH (tp)-[:tp_of_diagnosis]->(d)
WHERE tp.chemotherapy_regimen IN  ["FEC (3 week cycles)"]
WITH
        f, parent,p, ss, d,tp, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent,p, ss, d,tp, s, samp,
        f.file_size /(1024^i) AS value,
        10^precision AS factor,
        units[i] as unit
WITH
        f, parent,p, ss, d,tp, s, samp, unit,
        round(factor * value)/factor AS size
RETURN Distinct
    f.file_name AS `File Name`,
    head(labels(samp)) AS `Association`,
    f.file_description AS `Description`,
    f.file_format AS `File Format`,
     CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
    p.program_acronym AS `Program Code`,
    s.study_acronym AS `Arm`,
    ss.study_subject_id AS `Case ID`,
    samp.sample_id AS `Sample ID`
    order by f.file_name
'@

$statQuery = $ws.Cells.Item(2, 3).Value2
$neo4jFile = $ws.Cells.Item(2, 4).Value2
$webFile   = $ws.Cells.Item(2, 5).Value2

# Write the new tab-name labels first so they take the next shared-string
# slots, then fill in the rest of each row (matches the order the strings
# were appended to the workbook).
$ws.Cells.Item(3, 1).Value = "SamplesTab"
$ws.Cells.Item(4, 1).Value = "FilesTab"

# Row 3: SamplesTab
$ws.Cells.Item(3, 2).WrapText = $true
$ws.Cells.Item(3, 2).Value = $samplesQuery
$ws.Cells.Item(3, 3).WrapText = $true
$ws.Cells.Item(3, 3).Value = $statQuery
$ws.Cells.Item(3, 4).Value = $neo4jFile
$ws.Cells.Item(3, 5).Value = $webFile

# Row 4: FilesTab
$ws.Cells.Item(4, 2).WrapText = $true
$ws.Cells.Item(4, 2).Value = $filesQuery
$ws.Cells.Item(4, 3).WrapText = $true
$ws.Cells.Item(4, 3).Value = $statQuery
$ws.Cells.Item(4, 4).Value = $neo4jFile
$ws.Cells.Item(4, 5).Value = $webFile

# Match row heights to the ones Excel computed when it recalculated the
# sheet after adding these rows (wrapped-text row autosize).
$ws.Rows.Item(2).RowHeight = 316.8
$ws.Rows.Item(3).RowHeight = 331.2
$ws.Rows.Item(4).RowHeight = 409.6
